$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New week-30 header in column AG, matching the text-typed "1".."29" headers already
# present in row 1 (bold + centered, same as the rest of the header row).
$ws.Range("AF1").Copy($ws.Range("AG1"))
$ws.Range("AG1").NumberFormat = "@"
$ws.Range("AG1").Value = "30"
$ws.Range("AG1").NumberFormat = "General"

# Week-30 numeric counts for every facility row that already had a week-29 (AF) value.
$ws.Range("AG2").Value2 = 60
$ws.Range("AG4").Value2 = 1
$ws.Range("AG5").Value2 = 3
$ws.Range("AG6").Value2 = 117
$ws.Range("AG7").Value2 = 35
$ws.Range("AG8").Value2 = 38
$ws.Range("AG9").Value2 = 1
$ws.Range("AG10").Value2 = 4
$ws.Range("AG12").Value2 = 1
$ws.Range("AG14").Value2 = 2
$ws.Range("AG16").Value2 = 1
$ws.Range("AG17").Value2 = 2
$ws.Range("AG22").Value2 = 1
$ws.Range("AG23").Value2 = 2
$ws.Range("AG25").Value2 = 65
$ws.Range("AG28").Value2 = 15
$ws.Range("AG29").Value2 = 0
$ws.Range("AG30").Value2 = 23
$ws.Range("AG31").Value2 = 4
$ws.Range("AG32").Value2 = 8
$ws.Range("AG34").Value2 = 3
$ws.Range("AG35").Value2 = 39
$ws.Range("AG36").Value2 = 1
$ws.Range("AG37").Value2 = 10
$ws.Range("AG38").Value2 = 92
$ws.Range("AG39").Value2 = 1
$ws.Range("AG40").Value2 = 7
$ws.Range("AG41").Value2 = 9
$ws.Range("AG42").Value2 = 29
$ws.Range("AG43").Value2 = 196
$ws.Range("AG44").Value2 = 100
$ws.Range("AG45").Value2 = 172
$ws.Range("AG46").Value2 = 6
$ws.Range("AG47").Value2 = 104
$ws.Range("AG48").Value2 = 2
$ws.Range("AG49").Value2 = 0
$ws.Range("AG50").Value2 = 2
$ws.Range("AG52").Value2 = 31
$ws.Range("AG53").Value2 = 0
$ws.Range("AG54").Value2 = 0
$ws.Range("AG55").Value2 = 11
$ws.Range("AG56").Value2 = 55
$ws.Range("AG57").Value2 = 37
